$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Populate new "Table and columns" mapping detail (columns C/D/E) ---
# and update the renamed rows in column A (Extra-Ingredienten block, totals block).
# Writes are ordered so that brand-new shared-string values are first introduced
# in the same relative order as the target workbook, since the xlsx writer appends
# newly-seen string values to the shared-string table in first-use order.

$ws.Range("C3").Value = "Name"
$ws.Range("C4").Value = "PropertyKey"
$ws.Range("D4").Value = "PropertyValue"
$ws.Range("C5").Value = "Description"
$ws.Range("C14").Value = "DeliverySurcharge"
$ws.Range("D30").Value = "Restaurant"
$ws.Range("D31").Value = "Customer"
$ws.Range("C32").Value = "Phone"
$ws.Range("C35").Value = "City"
$ws.Range("D35").Value = "Address"
$ws.Range("C36").Value = "PlacementDate"
$ws.Range("D36").Value = "Order"
$ws.Range("C38").Value = "DeliveryDate"
$ws.Range("C40").Value = "ProductId"
$ws.Range("D40").Value = "OrderItem"
$ws.Range("C44").Value = "Amount"
$ws.Range("D45").Value = "OrderItemProperty"
$ws.Range("E45").Value = "note: OrderItemId van Product"
$ws.Range("C49").Value = "CouponId"
$ws.Range("A45").Value = "Extra-IngrediÃ«nten"
$ws.Range("A46").Value = "Prijs-Extra-IngrediÃ«nten"
$ws.Range("D3").Value = "Product"
$ws.Range("D5").Value = "Product"
$ws.Range("C6").Value = "Price"
$ws.Range("D6").Value = "Product"
$ws.Range("C9").Value = "Name"
$ws.Range("D9").Value = "Category"
$ws.Range("C10").Value = "Name"
$ws.Range("D10").Value = "Category"
$ws.Range("C11").Value = "Name"
$ws.Range("D11").Value = "Product"
$ws.Range("C12").Value = "Description"
$ws.Range("D12").Value = "Product"
$ws.Range("C13").Value = "Price"
$ws.Range("D13").Value = "Product"
$ws.Range("D14").Value = "Product"
$ws.Range("C15").Value = "PropertyKey"
$ws.Range("D15").Value = "PropertyValue"
$ws.Range("C16").Value = "PropertyKey"
$ws.Range("D16").Value = "PropertyValue"
$ws.Range("C17").Value = "PropertyKey"
$ws.Range("D17").Value = "PropertyValue"
$ws.Range("C18").Value = "PropertyKey"
$ws.Range("D18").Value = "PropertyValue"
$ws.Range("C21").Value = "Name"
$ws.Range("D21").Value = "Category"
$ws.Range("C22").Value = "Name"
$ws.Range("D22").Value = "Category"
$ws.Range("C23").Value = "Name"
$ws.Range("D23").Value = "Product"
$ws.Range("C24").Value = "Description"
$ws.Range("D24").Value = "Product"
$ws.Range("C25").Value = "Price"
$ws.Range("D25").Value = "Product"
$ws.Range("C26").Value = "PropertyKey"
$ws.Range("D26").Value = "PropertyValue"
$ws.Range("C27").Value = "PropertyKey"
$ws.Range("D27").Value = "PropertyValue"
$ws.Range("C30").Value = "Name"
$ws.Range("C31").Value = "Name"
$ws.Range("D32").Value = "Customer"
$ws.Range("C33").Value = "Email"
$ws.Range("D33").Value = "Customer"
$ws.Range("D38").Value = "Order"
$ws.Range("C39").Value = "DeliveryDate"
$ws.Range("D39").Value = "Order"
$ws.Range("C41").Value = "ProductId"
$ws.Range("D41").Value = "OrderItem"
$ws.Range("C42").Value = "ProductId"
$ws.Range("D42").Value = "OrderItem"
$ws.Range("D44").Value = "OrderItem"
$ws.Range("C45").Value = "PropertyKey"
$ws.Range("A47").Value = "Regelprijs"
$ws.Range("A48").Value = "Totaalprijs"
$ws.Range("A49").Value = "Gebruikte Coupon"
$ws.Range("D49").Value = "Order"
$ws.Range("A50").Value = "Coupon Korting"
$ws.Range("A51").Value = "Te Betalen"

# --- Remove the three now-obsolete rows (old Coupon/Korting/TeBetalen trailer) ---
# Use Clear() (not ClearContents) so the rows disappear from the sparse row list
# entirely, instead of surviving as empty-but-styled rows (matches target diff, which
# drops rows 52-54 outright while rows 55-57 keep their original row numbers).
$ws.Range("A52").Clear()
$ws.Range("A53").Clear()
$ws.Range("A54").Clear()

# --- Column widths for the two new columns (and minor width tweaks on A/B) ---
$ws.Columns.Item(1).ColumnWidth = 23.666666666666668
$ws.Columns.Item(2).ColumnWidth = 16.833333333333332
$ws.Columns.Item(3).ColumnWidth = 18.5
$ws.Columns.Item(4).ColumnWidth = 17.5

# --- Final selection, matching the saved workbook state ---
$ws.Range("E45").Select()
